$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "nityaranjn623@gmail.com"
$ws.Range("C3").Value = "abh090824@gmail.com"

$ws.Range("C6").Select()
